$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("326:326").Insert()

$ws.Cells.Item(326, 1).Value = 6
$ws.Cells.Item(326, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(326, 3).Value = "Metropolitana"
$ws.Cells.Item(326, 4).Value = 45194
$ws.Cells.Item(326, 5).Value = 13
$ws.Cells.Item(326, 6).Value = 100112001
$ws.Cells.Item(326, 7).Value = "Berenjena"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 155
$ws.Cells.Item(326, 11).Value = 6500
$ws.Cells.Item(326, 12).Value = 7000
$ws.Cells.Item(326, 13).Value = 6823
$ws.Cells.Item(326, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(326, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(326, 16).Value = 136
$ws.Cells.Item(326, 17).Value = 50
$ws.Cells.Item(326, 18).Value = "Hortaliza"
